# Add two explanatory note textboxes under two of the charts on slide 1,
# as described in the commit: "Added a note to the graphs we can
# potentially make interactive".
#
# AddTextbox's Left/Top/Width/Height parameters are in points, while the
# target offsets/extents from the OOXML are in EMU (914400 EMU = 1 inch,
# 12700 EMU = 1 point) - convert before calling.

function EmuToPt([double]$emu) {
    return $emu / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$noteText = "* Opportunity to make interactive: can allow users to select between " + [char]0x2018 + "good" + [char]0x2019 + " and " + [char]0x2018 + "bad" + [char]0x2019 + " category from dropdown to see stats for selection."

# --- Textbox 1: note under the "Average of balance scores overtime" chart ---
$tb1 = $s.Shapes.AddTextbox(
    1,
    (EmuToPt 8394086),
    (EmuToPt 3868118),
    (EmuToPt 3998277),
    (EmuToPt 400110)
)
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1
$tb1.Fill.Visible = $false
$tb1.TextFrame.TextRange.Text = $noteText
$tb1.TextFrame.TextRange.Font.Size = 10
$tb1.TextFrame.TextRange.LanguageID = "en-CA"

# --- Textbox 2: note under the "Frequency of categorical responses..." chart ---
$tb2 = $s.Shapes.AddTextbox(
    1,
    (EmuToPt 2623221),
    (EmuToPt 4746289),
    (EmuToPt 3760892),
    (EmuToPt 400110)
)
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1
$tb2.Fill.Visible = $false
$tb2.TextFrame.TextRange.Text = $noteText
$tb2.TextFrame.TextRange.Font.Size = 10
$tb2.TextFrame.TextRange.LanguageID = "en-CA"

Write-Output "Added shapes: $($tb1.Id)=$($tb1.Name), $($tb2.Id)=$($tb2.Name)"
